$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new columns before the existing "bir" column (C), shifting
# bir/12r (and everything after) two columns to the right.
$ws.Range("C1:D1").EntireColumn.Insert()

# New column headers: 12r-adj (C) and bir-adj (D) -- copies of the
# (now shifted) 12r/bir columns, added as new "adjusted" response columns.
$ws.Range("C1").Value = "12r-adj"
$ws.Range("D1").Value = "bir-adj"

# Populate the new columns with the original bir/12r values (now living
# in columns E and F respectively), swapped in order: 12r-adj takes the
# 12r values, bir-adj takes the bir values.
for ($r = 2; $r -le 7; $r++) {
    $birVal = $ws.Cells($r, 5).Value()
    $r12Val = $ws.Cells($r, 6).Value()
    $ws.Cells($r, 3).Value = $r12Val
    $ws.Cells($r, 4).Value = $birVal
}

# Match the column widths the author set on Austin-Num / 12r-adj / bir-adj.
$ws.Range("B1:D1").EntireColumn.ColumnWidth = 9.25

# Leave the selection where the author left it after typing the last value.
$ws.Range("E8").Select()
